$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2538.1667
$ws.Range("I64").Value = 2401.9
$ws.Range("J64").Value = 2681.6052
$ws.Range("K64").Value = 2401.9
$ws.Range("L64").Value = 2681.6052
$ws.Range("M64").Value = -2153.9
$ws.Range("N64").Value = -3177.6052
$ws.Range("H67").Value = 2538.1667
$ws.Range("I67").Value = 2401.9
$ws.Range("J67").Value = 2681.6052
$ws.Range("K67").Value = 2401.9
$ws.Range("L67").Value = 2681.6052
$ws.Range("M67").Value = -1543.9
$ws.Range("N67").Value = -4397.6052
$ws.Range("H74").Value = 3377.6667
$ws.Range("I74").Value = 2914.1428
$ws.Range("J74").Value = 5000
$ws.Range("K74").Value = 2914.1428
$ws.Range("L74").Value = 5000
$ws.Range("M74").Value = -1978.1428
$ws.Range("N74").Value = -6872
$ws.Range("H77").Value = 3377.6667
$ws.Range("I77").Value = 2914.1428
$ws.Range("J77").Value = 5000
$ws.Range("K77").Value = 14570.714
$ws.Range("L77").Value = 25000
$ws.Range("M77").Value = -9890.714
$ws.Range("N77").Value = -34360
$ws.Range("H80").Value = 376.47058
$ws.Range("I80").Value = 340.06668
$ws.Range("J80").Value = 649.5
$ws.Range("K80").Value = 1020.20004
$ws.Range("L80").Value = 1948.5
$ws.Range("M80").Value = -22.20004000000006
$ws.Range("N80").Value = -3944.5
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H83").Value = 376.47058
$ws.Range("I83").Value = 340.06668
$ws.Range("J83").Value = 649.5
$ws.Range("K83").Value = 3060.60012
$ws.Range("L83").Value = 5845.5
$ws.Range("M83").Value = 1931.39988
$ws.Range("N83").Value = -15829.5
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H103").Value = 886.6667
$ws.Range("I103").Value = 1500
$ws.Range("J103").Value = 580
$ws.Range("K103").Value = 4500
$ws.Range("L103").Value = 1740
$ws.Range("M103").Value = -3914
$ws.Range("N103").Value = -2912
$ws.Range("H138").Value = 2091.8108
$ws.Range("I138").Value = 1097.0476
$ws.Range("J138").Value = 3397.4375
$ws.Range("K138").Value = 3291.142800000001
$ws.Range("L138").Value = 10192.3125
$ws.Range("M138").Value = 1848.857199999999
$ws.Range("N138").Value = -20472.3125
$ws.Range("H141").Value = 3214.6562
$ws.Range("I141").Value = 1743.1111
$ws.Range("J141").Value = 11161
$ws.Range("K141").Value = 5229.3333
$ws.Range("L141").Value = 33483
$ws.Range("M141").Value = -49.33330000000024
$ws.Range("N141").Value = -43843

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1293.8064
$ws.Range("I61").Value = 1392.625
$ws.Range("J61").Value = 1188.4
$ws.Range("K61").Value = 1392.625
$ws.Range("L61").Value = 1188.4
$ws.Range("M61").Value = -1180.625
$ws.Range("N61").Value = -1612.4
$ws.Range("H93").Value = 36000
$ws.Range("J93").Value = 36000
$ws.Range("L93").Value = 36000
$ws.Range("N93").Value = -40992
$ws.Range("H97").Value = 924.67566
$ws.Range("I97").Value = 719.7406999999999
$ws.Range("J97").Value = 1478
$ws.Range("K97").Value = 719.7406999999999
$ws.Range("L97").Value = 1478
$ws.Range("M97").Value = -223.7406999999999
$ws.Range("N97").Value = -2470
$ws.Range("H132").Value = 19609678
$ws.Range("I132").Value = 27028630
$ws.Range("J132").Value = 3462549.5
$ws.Range("K132").Value = 81085890
$ws.Range("L132").Value = 10387648.5
$ws.Range("M132").Value = -81083360
$ws.Range("N132").Value = -10392708.5
$ws.Range("H136").Value = 1293.8064
$ws.Range("I136").Value = 1392.625
$ws.Range("J136").Value = 1188.4
$ws.Range("K136").Value = 4177.875
$ws.Range("L136").Value = 3565.2
$ws.Range("M136").Value = -1627.875
$ws.Range("N136").Value = -8665.200000000001
$ws.Range("H9").Value = 60000
$ws.Range("J9").Value = 60000
$ws.Range("L9").Value = 60000
$ws.Range("N9").Value = -60336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1017.0909
$ws.Range("I94").Value = 1120
$ws.Range("J94").Value = 931.3333
$ws.Range("K94").Value = 1120
$ws.Range("L94").Value = 931.3333
$ws.Range("M94").Value = -669
$ws.Range("N94").Value = -1833.3333
$ws.Range("H105").Value = 76924960
$ws.Range("I105").Value = 1869
$ws.Range("K105").Value = 1869
$ws.Range("M105").Value = -122
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("H134").Value = 2367180
$ws.Range("I134").Value = 810.19354
$ws.Range("J134").Value = 6952021
$ws.Range("K134").Value = 2430.58062
$ws.Range("L134").Value = 20856063
$ws.Range("M134").Value = 104.4193800000003
$ws.Range("N134").Value = -20861133
$ws.Range("M113").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1357.2439
$ws.Range("I31").Value = 988.7917
$ws.Range("J31").Value = 1877.4117
$ws.Range("K31").Value = 988.7917
$ws.Range("L31").Value = 1877.4117
$ws.Range("M31").Value = -693.7917
$ws.Range("N31").Value = -2467.4117
$ws.Range("H34").Value = 1357.2439
$ws.Range("I34").Value = 988.7917
$ws.Range("J34").Value = 1877.4117
$ws.Range("K34").Value = 988.7917
$ws.Range("L34").Value = 1877.4117
$ws.Range("M34").Value = -786.7917
$ws.Range("N34").Value = -2281.4117
$ws.Range("H74").Value = 43314
$ws.Range("J74").Value = 43314
$ws.Range("L74").Value = 43314
$ws.Range("N74").Value = -45062
$ws.Range("H77").Value = 43314
$ws.Range("J77").Value = 43314
$ws.Range("L77").Value = 129942
$ws.Range("N77").Value = -138678

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 3904.9167
$ws.Range("J103").Value = 4638.4
$ws.Range("L103").Value = 13915.2
$ws.Range("N103").Value = -15673.2
$ws.Range("H113").Value = 34444884
$ws.Range("I113").Value = 27778048
$ws.Range("J113").Value = 36111596
$ws.Range("K113").Value = 83334144
$ws.Range("L113").Value = 108334788
$ws.Range("M113").Value = -83331974
$ws.Range("N113").Value = -108339128

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 9052
$ws.Range("I99").Value = 1382.8572
$ws.Range("J99").Value = 17999.334
$ws.Range("K99").Value = 1382.8572
$ws.Range("L99").Value = 17999.334
$ws.Range("M99").Value = 863.1428000000001
$ws.Range("N99").Value = -22491.334
$ws.Range("H132").Value = 7835.143
$ws.Range("I132").Value = 4893.0835
$ws.Range("J132").Value = 25487.5
$ws.Range("K132").Value = 14679.2505
$ws.Range("L132").Value = 76462.5
$ws.Range("M132").Value = -12149.2505
$ws.Range("N132").Value = -81522.5
$ws.Range("H137").Value = 39340
$ws.Range("I137").Value = 29680
$ws.Range("J137").Value = 49000
$ws.Range("K137").Value = 29680
$ws.Range("L137").Value = 49000
$ws.Range("M137").Value = -24580
$ws.Range("N137").Value = -59200

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("H74").Value = 15333.333
$ws.Range("J74").Value = 18000
$ws.Range("L74").Value = 18000
$ws.Range("N74").Value = -19996
$ws.Range("H77").Value = 15333.333
$ws.Range("J77").Value = 18000
$ws.Range("L77").Value = 54000
$ws.Range("N77").Value = -63984
$ws.Range("H80").Value = 49800
$ws.Range("J80").Value = 49800
$ws.Range("L80").Value = 49800
$ws.Range("N80").Value = -52046
$ws.Range("H83").Value = 49800
$ws.Range("J83").Value = 49800
$ws.Range("L83").Value = 149400
$ws.Range("N83").Value = -160632
$ws.Range("H93").Value = 1450
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 1450
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 1450
$ws.Range("N93").Value = -3946
$ws.Range("H109").Value = 18493.375
$ws.Range("J109").Value = 18493.375
$ws.Range("L109").Value = 18493.375
$ws.Range("N109").Value = -21267.375
$ws.Range("H127").Value = 52075
$ws.Range("J127").Value = 52075
$ws.Range("L127").Value = 52075
$ws.Range("N127").Value = -61995
$ws.Range("N21").ClearContents()
$ws.Range("M93").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 17746.479
$ws.Range("J14").Value = 17746.479
$ws.Range("L14").Value = 17746.479
$ws.Range("N14").Value = -18082.479
$ws.Range("H108").Value = 36647.2
$ws.Range("J108").Value = 36647.2
$ws.Range("L108").Value = 36647.2
$ws.Range("N108").Value = -44327.2
$ws.Range("H132").Value = 19467.322
$ws.Range("I132").Value = 26527.1
$ws.Range("J132").Value = 6631.364
$ws.Range("K132").Value = 79581.29999999999
$ws.Range("L132").Value = 19894.092
$ws.Range("M132").Value = -77051.29999999999
$ws.Range("N132").Value = -24954.092
